$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.462.14'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '1.647.89'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '299.67'
$ws.Range("E6").Value = '  -1.58%  '
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3559'
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08100'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.219'
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.00'
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.397'
$ws.Range("E14").Value = '  -2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.376'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001194'
$ws.Range("E16").Value = '  -3.03%  '
$ws.Range("D17").Value = '1.652.69'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.26'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06969'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.738'
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.41'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").Value = '23.488.62'
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.908'
$ws.Range("E26").Value = '  -5.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.91'
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.62'
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.202'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.84'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").Value = '1.834.14'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.902'
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.89'
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.109'
$ws.Range("E34").Value = '  +0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.015'
$ws.Range("E35").Value = '  -6.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02723'
$ws.Range("E36").Value = '  -2.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08739'
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2435'
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.942'
$ws.Range("E39").Value = '  -2.95%  '
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06791'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6889'
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.316'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.54'
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6399'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.261'
$ws.Range("E47").Value = '  -3.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.921'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07737'
$ws.Range("E49").Value = '  -3.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.52'
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.147'
$ws.Range("E51").Value = '  -3.77%  '
